$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The "Prüfungsform" (column S) marker ("Klausur") is no longer needed for the
#    data rows - remove those cells outright (content + formatting), which also
#    drops the now-unused "Klausur" string from the shared strings table.
$ws.Range("S6:S11").Clear()

# 2. Update Matrikelnummer (column F) for the re-sorted / corrected students.
$ws.Range("F7").Value = 300002
$ws.Range("F8").Value = 300003
$ws.Range("F9").Value = 300004
$ws.Range("F10").Value = 310005
$ws.Range("F11").Value = 300013

# 3. Enter the grades (column G, "Leistung") for the students that now have one.
$ws.Range("G7").Value = 3.7
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 1
$ws.Range("G11").Value = 5

# The cells that now hold a hand-entered Matrikelnummer switch from the
# shaded/locked look to the plain unlocked look already used by column G.
# Copy the formatting (not the value) from the matching G cell in each row.
$ws.Range("G7").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("G11").Copy()
$ws.Range("F11").PasteSpecial(-4122)

# 4. A new (empty) row further down the sheet, formatted like the rest of column F/G.
$ws.Range("G6").Copy()
$ws.Range("F25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# 5. Move the active selection to match the author's last cursor position.
$ws.Range("G9").Select()
